$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column M: "Voltage Divider Low" header + per-row divider formula.
$ws.Range("M1").Value = "Voltage Divider Low"
$ws.Range("M2").Formula = "=3.3*B2/(B2+10000)"
$ws.Range("M3:M17").Formula = "=3.3*B3/(B3+10000)"

# View changes: zoom level, scroll back to top-left (drops topLeftCell), and
# move the active selection to M20.
$excel.ActiveWindow.Zoom = 74
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("M20").Select() | Out-Null
